$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I (municipio-nombre) dimension/measure info
$ws.Range("I2").Value = "sdmx-dimension:refArea"
$ws.Range("I3").Value = "dim"
$ws.Range("I4").Value = "URI-Municipio"

# Column K (sexo) dimension/measure info
$ws.Range("K2").Value = "iaest-measure:sexo"
$ws.Range("K3").Value = "medida"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("K5").Clear()

# Column L (direccion-provincial-nombre) dimension/measure info
$ws.Range("L2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("L3").Value = "medida"
$ws.Range("L4").Value = "xsd:int"
